$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 47836.383
$ws.Range("I33").Value = 59077.707
$ws.Range("J33").Value = 60.75
$ws.Range("K33").Value = 59077.707
$ws.Range("L33").Value = 60.75
$ws.Range("M33").Value = -58848.707
$ws.Range("N33").Value = -518.75

$ws.Range("H64").Value = 6215.8125
$ws.Range("J64").Value = 5034.846
$ws.Range("L64").Value = 5034.846
$ws.Range("N64").Value = -5530.846

$ws.Range("H67").Value = 6215.8125
$ws.Range("J67").Value = 5034.846
$ws.Range("L67").Value = 5034.846
$ws.Range("N67").Value = -6750.846

$ws.Range("H74").Value = 4723638.5
$ws.Range("I74").Value = 8654754
$ws.Range("J74").Value = 6300
$ws.Range("K74").Value = 8654754
$ws.Range("L74").Value = 6300
$ws.Range("M74").Value = -8653818
$ws.Range("N74").Value = -8172

$ws.Range("H76").Value = 157148720
$ws.Range("I76").Value = 157148720
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 157148720
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -157148405

$ws.Range("H77").Value = 4723638.5
$ws.Range("I77").Value = 8654754
$ws.Range("J77").Value = 6300
$ws.Range("K77").Value = 43273770
$ws.Range("L77").Value = 31500
$ws.Range("M77").Value = -43269090
$ws.Range("N77").Value = -40860

$ws.Range("H79").Value = 157148720
$ws.Range("I79").Value = 157148720
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 157148720
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -157147628

$ws.Range("H88").Value = 722386.3
$ws.Range("I88").Value = 1014.7143
$ws.Range("J88").Value = 1059026.4
$ws.Range("K88").Value = 1014.7143
$ws.Range("L88").Value = 1059026.4
$ws.Range("M88").Value = -608.7143
$ws.Range("N88").Value = -1059838.4

$ws.Range("H91").Value = 722386.3
$ws.Range("I91").Value = 1014.7143
$ws.Range("J91").Value = 1059026.4
$ws.Range("K91").Value = 1014.7143
$ws.Range("L91").Value = 1059026.4
$ws.Range("M91").Value = 389.2857
$ws.Range("N91").Value = -1061834.4

$ws.Range("H116").Value = 2515.0908
$ws.Range("I116").Value = 2432.5
$ws.Range("J116").Value = 2735.3333
$ws.Range("K116").Value = 2432.5
$ws.Range("L116").Value = 2735.3333
$ws.Range("M116").Value = 1009.5
$ws.Range("N116").Value = -9619.3333

$ws.Range("H129").Value = 1317.7878
$ws.Range("I129").Value = 346
$ws.Range("J129").Value = 1380.4839
$ws.Range("K129").Value = 1038
$ws.Range("L129").Value = 4141.4517
$ws.Range("M129").Value = 3962
$ws.Range("N129").Value = -14141.4517

$ws.Range("H138").Value = 3330.79
$ws.Range("I138").Value = 823.55554
$ws.Range("J138").Value = 3881.1584
$ws.Range("K138").Value = 2470.66662
$ws.Range("L138").Value = 11643.4752
$ws.Range("M138").Value = 2669.33338
$ws.Range("N138").Value = -21923.4752

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1895.9474
$ws.Range("I61").Value = 2081.6
$ws.Range("J61").Value = 1199.75
$ws.Range("K61").Value = 2081.6
$ws.Range("L61").Value = 1199.75
$ws.Range("M61").Value = -1869.6
$ws.Range("N61").Value = -1623.75

$ws.Range("H136").Value = 1895.9474
$ws.Range("I136").Value = 2081.6
$ws.Range("J136").Value = 1199.75
$ws.Range("K136").Value = 6244.799999999999
$ws.Range("L136").Value = 3599.25
$ws.Range("M136").Value = -3694.799999999999
$ws.Range("N136").Value = -8699.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 472.5
$ws.Range("I22").Value = 383.75
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 383.75
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -210.75
$ws.Range("N22").Value = -996

$ws.Range("H80").Value = 698
$ws.Range("I80").Value = 400.85715
$ws.Range("J80").Value = 958
$ws.Range("K80").Value = 400.85715
$ws.Range("L80").Value = 958
$ws.Range("M80").Value = 597.14285
$ws.Range("N80").Value = -2954

$ws.Range("H83").Value = 698
$ws.Range("I83").Value = 400.85715
$ws.Range("J83").Value = 958
$ws.Range("K83").Value = 2004.28575
$ws.Range("L83").Value = 4790
$ws.Range("M83").Value = 2987.71425
$ws.Range("N83").Value = -14774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 238.3
$ws.Range("I22").Value = 252.92308
$ws.Range("J22").Value = 211.14285
$ws.Range("K22").Value = 252.92308
$ws.Range("L22").Value = 211.14285
$ws.Range("M22").Value = 97.07692
$ws.Range("N22").Value = -911.14285

$ws.Range("H86").Value = 3988
$ws.Range("I86").Value = 3984.5715
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 3984.5715
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -2861.5715
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 3988
$ws.Range("I89").Value = 3984.5715
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 19922.8575
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -14306.8575
$ws.Range("N89").Value = -31232

$ws.Range("H99").Value = 1916.925
$ws.Range("I99").Value = 1547.9166
$ws.Range("J99").Value = 2470.4375
$ws.Range("K99").Value = 1547.9166
$ws.Range("L99").Value = 2470.4375
$ws.Range("M99").Value = -49.91660000000002
$ws.Range("N99").Value = -5466.4375

$ws.Range("H126").Value = 1916.925
$ws.Range("I126").Value = 1547.9166
$ws.Range("J126").Value = 2470.4375
$ws.Range("K126").Value = 4643.7498
$ws.Range("L126").Value = 7411.3125
$ws.Range("M126").Value = -2173.7498
$ws.Range("N126").Value = -12351.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1263.1305
$ws.Range("J5").Value = 4733.3335
$ws.Range("L5").Value = 14200.0005
$ws.Range("N5").Value = -14424.0005

$ws.Range("H14").Value = 104.15385
$ws.Range("I14").Value = 104.15385
$ws.Range("K14").Value = 312.46155
$ws.Range("M14").Value = -139.46155

$ws.Range("H132").Value = 455213.88
$ws.Range("I132").Value = 606.25
$ws.Range("J132").Value = 1667500.9
$ws.Range("K132").Value = 5456.25
$ws.Range("L132").Value = 15007508.1
$ws.Range("M132").Value = -2926.25
$ws.Range("N132").Value = -15012568.1

$ws.Range("H135").Value = 1263.1305
$ws.Range("J135").Value = 4733.3335
$ws.Range("L135").Value = 42600.0015
$ws.Range("N135").Value = -47670.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 919583.8
$ws.Range("I107").Value = 1634426.5
$ws.Range("J107").Value = 500.2857
$ws.Range("K107").Value = 1634426.5
$ws.Range("L107").Value = 500.2857
$ws.Range("M107").Value = -1632506.5
$ws.Range("N107").Value = -4340.2857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1664.2142
$ws.Range("I7").Value = 1676.8462
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1676.8462
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -1564.8462
$ws.Range("N7").Value = -1724

$ws.Range("H22").Value = 637.61536
$ws.Range("I22").Value = 653.63635
$ws.Range("J22").Value = 549.5
$ws.Range("K22").Value = 653.63635
$ws.Range("L22").Value = 549.5
$ws.Range("M22").Value = -358.63635
$ws.Range("N22").Value = -1139.5

$ws.Range("H27").Value = 637.61536
$ws.Range("I27").Value = 653.63635
$ws.Range("J27").Value = 549.5
$ws.Range("K27").Value = 653.63635
$ws.Range("L27").Value = 549.5
$ws.Range("M27").Value = -546.63635
$ws.Range("N27").Value = -763.5

$ws.Range("H126").Value = 1664.2142
$ws.Range("I126").Value = 1676.8462
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 5030.5386
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -2560.5386
$ws.Range("N126").Value = -9440

$ws.Range("H132").Value = 2812.5715
$ws.Range("I132").Value = 2876.75
$ws.Range("J132").Value = 2672.5454
$ws.Range("K132").Value = 8630.25
$ws.Range("L132").Value = 8017.6362
$ws.Range("M132").Value = -6100.25
$ws.Range("N132").Value = -13077.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1206.4286
$ws.Range("I132").Value = 838.13043
$ws.Range("J132").Value = 2900.6
$ws.Range("K132").Value = 2514.39129
$ws.Range("L132").Value = 8701.799999999999
$ws.Range("M132").Value = 15.60870999999997
$ws.Range("N132").Value = -13761.8
